$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 28 (current row 28 -> becomes row 30 after insert)
$ws.Range("A28:A29").EntireRow.Insert()

# New row 28 data
$ws.Cells.Item(28, 1).Value = 4
$ws.Cells.Item(28, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(28, 3).Value = "Los Lagos"
$ws.Cells.Item(28, 4).Value = 44897
$ws.Cells.Item(28, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(28, 5).Value = 10
$ws.Cells.Item(28, 6).Value = 300000000
$ws.Cells.Item(28, 7).Value = "Espárragos"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 300
$ws.Cells.Item(28, 11).Value = 2000
$ws.Cells.Item(28, 12).Value = 2000
$ws.Cells.Item(28, 13).Value = 2000
$ws.Cells.Item(28, 14).Value = "`$/kilo"
$ws.Cells.Item(28, 15).Value = "Provincia de Linares"
$ws.Cells.Item(28, 16).Value = 2000
$ws.Cells.Item(28, 17).Value = 1
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# New row 29 data
$ws.Cells.Item(29, 1).Value = 4
$ws.Cells.Item(29, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(29, 3).Value = "Los Lagos"
$ws.Cells.Item(29, 4).Value = 44897
$ws.Cells.Item(29, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 5).Value = 10
$ws.Cells.Item(29, 6).Value = 300000000
$ws.Cells.Item(29, 7).Value = "Espárragos"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 300
$ws.Cells.Item(29, 11).Value = 1500
$ws.Cells.Item(29, 12).Value = 1500
$ws.Cells.Item(29, 13).Value = 1500
$ws.Cells.Item(29, 14).Value = "`$/kilo"
$ws.Cells.Item(29, 15).Value = "Provincia de Linares"
$ws.Cells.Item(29, 16).Value = 1500
$ws.Cells.Item(29, 17).Value = 1
$ws.Cells.Item(29, 18).Value = "Hortaliza"
